$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain numeric strings (e.g. "5.05", "1.00").
# Force these specific cells to Text format first so Excel keeps the literal
# string (incl. trailing zeros) instead of silently converting to a Number.
$textForceCells = @(
    'D5',
    'D6',
    'D10',
    'D12',
    'D13',
    'D14',
    'D16',
    'D19',
    'D20',
    'D21',
    'D22',
    'D23',
    'D25',
    'D26',
    'D27',
    'D28',
    'D31',
    'D32',
    'D33',
    'D34',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D47',
    'D48',
    'D49',
    'D50',
    'D51'
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '60.891.16'
$ws.Range('E2').Value = '  -1.71%  '

$ws.Range('D3').Value = '2.413.87'
$ws.Range('E3').Value = '  -1.23%  '

$ws.Range('D5').Value = '567.50'
$ws.Range('E5').Value = '  -2.24%  '

$ws.Range('D6').Value = '138.98'
$ws.Range('E6').Value = '  -2.18%  '

$ws.Range('E7').Value = '  +0.18%  '

$ws.Range('E8').Value = '  -0.45%  '

$ws.Range('D9').Value = '2.397.64'
$ws.Range('E9').Value = '  -1.66%  '

$ws.Range('D10').Value = '0.106'
$ws.Range('E10').Value = '  -3.70%  '

$ws.Range('E11').Value = '  -0.42%  '

$ws.Range('D12').Value = '5.05'
$ws.Range('E12').Value = '  -2.90%  '

$ws.Range('D13').Value = '0.335'
$ws.Range('E13').Value = '  -1.59%  '

$ws.Range('D14').Value = '25.98'
$ws.Range('E14').Value = '  -0.97%  '

$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.862.32'
$ws.Range('E15').Value = '  -1.05%  '

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.0000170'
$ws.Range('E16').Value = '  -1.99%  '

$ws.Range('D17').Value = '60.872.74'
$ws.Range('E17').Value = '  -1.81%  '

$ws.Range('D18').Value = '2.307.08'
$ws.Range('E18').Value = '  -5.33%  '

$ws.Range('D19').Value = '7.89'
$ws.Range('E19').Value = '  +9.64%  '

$ws.Range('D20').Value = '10.53'
$ws.Range('E20').Value = '  -0.93%  '

$ws.Range('D21').Value = '321.86'
$ws.Range('E21').Value = '  -0.96%  '

$ws.Range('D22').Value = '4.03'
$ws.Range('E22').Value = '  -0.97%  '

$ws.Range('D23').Value = '6.17'
$ws.Range('E23').Value = '  +3.24%  '

$ws.Range('E24').Value = '  +0.09%  '

$ws.Range('D25').Value = '1.82'
$ws.Range('E25').Value = '  -3.75%  '

$ws.Range('D26').Value = '64.18'
$ws.Range('E26').Value = '  -1.85%  '

$ws.Range('D27').Value = '580.89'
$ws.Range('E27').Value = '  -1.32%  '

$ws.Range('D28').Value = '8.25'
$ws.Range('E28').Value = '  -9.43%  '

$ws.Range('D29').Value = '2.544.97'
$ws.Range('E29').Value = '  -0.77%  '

$ws.Range('D30').Value = '0.0₃0921'
$ws.Range('E30').Value = '  -3.09%  '

$ws.Range('D31').Value = '7.87'
$ws.Range('E31').Value = '  -0.68%  '

$ws.Range('D32').Value = '1.33'
$ws.Range('E32').Value = '  -4.57%  '

$ws.Range('D33').Value = '1.81'
$ws.Range('E33').Value = '  -3.99%  '

$ws.Range('D34').Value = '0.132'
$ws.Range('E34').Value = '  -1.63%  '

$ws.Range('E35').Value = '  -0.06%  '

$ws.Range('D36').Value = '1.41'
$ws.Range('E36').Value = '  +0.11%  '

$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '151.29'
$ws.Range('E37').Value = '  -0.96%  '

$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '4.57'
$ws.Range('E38').Value = '  -4.87%  '

$ws.Range('D39').Value = '0.366'
$ws.Range('E39').Value = '  -2.04%  '

$ws.Range('D40').Value = '18.14'
$ws.Range('E40').Value = '  -1.21%  '

$ws.Range('D41').Value = '5.11'
$ws.Range('E41').Value = '  -2.04%  '

$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.06%  '

$ws.Range('D43').Value = '1.66'
$ws.Range('E43').Value = '  -1.69%  '

$ws.Range('D44').Value = '41.13'
$ws.Range('E44').Value = '  -4.43%  '

$ws.Range('D45').Value = '2.38'
$ws.Range('E45').Value = '  -3.00%  '

$ws.Range('D46').Value = '0.0₆0294'
$ws.Range('E46').Value = '  +8.02%  '

$ws.Range('D47').Value = '142.39'
$ws.Range('E47').Value = '  +1.00%  '

$ws.Range('D48').Value = '3.51'
$ws.Range('E48').Value = '  -2.32%  '

$ws.Range('D49').Value = '0.584'
$ws.Range('E49').Value = '  -2.49%  '

$ws.Range('D50').Value = '19.34'
$ws.Range('E50').Value = '  -1.46%  '

$ws.Range('D51').Value = '0.0498'
$ws.Range('E51').Value = '  -3.26%  '
